$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Velocity + Weight")

# New SKU rows appended to the bottom of the "Velocity + Weight" sheet
# (rows 567-591). Columns populated: A (SKU id), B (velocity count, manually
# entered, plain black font), C (z-scaled velocity formula), G (C + 0.5306),
# H (literal 0 - these new SKUs have no weight-side data yet).
# "BRef" is the row number actually referenced inside the C-column formula -
# it matches the row itself for every new record except row 572, whose
# formula was typed/filled referencing B72 instead of B572.
$newRows = @(
    @{Row=567; Sku="CS56603";      B=97;  BRef=567},
    @{Row=568; Sku="EA59121";      B=82;  BRef=568},
    @{Row=569; Sku="CS58103";      B=48;  BRef=569},
    @{Row=570; Sku="EA56410";      B=58;  BRef=570},
    @{Row=571; Sku="CS61400";      B=41;  BRef=571},
    @{Row=572; Sku="EA56103";      B=18;  BRef=72},
    @{Row=573; Sku="CS56103";      B=33;  BRef=573},
    @{Row=574; Sku="EA58510";      B=73;  BRef=574},
    @{Row=575; Sku="CS57603";      B=58;  BRef=575},
    @{Row=576; Sku="EA57610";      B=134; BRef=576},
    @{Row=577; Sku="CS61300";      B=29;  BRef=577},
    @{Row=578; Sku="PKXX40#";      B=0;   BRef=578},
    @{Row=579; Sku="CS58503";      B=16;  BRef=579},
    @{Row=580; Sku="PK39226BAGr2"; B=0;   BRef=580},
    @{Row=581; Sku="CS59503";      B=21;  BRef=581},
    @{Row=582; Sku="CS57903";      B=74;  BRef=582},
    @{Row=583; Sku="CS56403";      B=47;  BRef=583},
    @{Row=584; Sku="EA59110";      B=85;  BRef=584},
    @{Row=585; Sku="CS57103";      B=58;  BRef=585},
    @{Row=586; Sku="EA57103";      B=23;  BRef=586},
    @{Row=587; Sku="EA57110";      B=87;  BRef=587},
    @{Row=588; Sku="CS56503";      B=68;  BRef=588},
    @{Row=589; Sku="EA59510";      B=76;  BRef=589},
    @{Row=590; Sku="CS59103";      B=25;  BRef=590},
    @{Row=591; Sku="CS61100";      B=41;  BRef=591}
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Range("A$r").Value = $item.Sku

    $bCell = $ws.Range("B$r")
    $bCell.Value = $item.B
    $bCell.Font.Color = 0

    $ws.Range("C$r").Formula = "=((B$($item.BRef) - L2)/L3)"
    $ws.Range("G$r").Formula = "=C$r+0.5306"
    $ws.Range("H$r").Value = 0
}

# Reflect the cursor position left behind after the data entry.
$ws.Range("G599").Select() | Out-Null
